$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keywords")

# --- Row 4: height change, and D4 content (new multi-line keywords text) ---
$ws.Rows.Item(4).RowHeight = 60
$ws.Range("D4").Value = "USE mydabases;`nSHOW TABLES`nSHOW COLUMNS FROM mytable`nDESC mytable # DESC here means describe"

# --- Row 6: height change, and D6 content (new CREATE TABLE example with USE) ---
$ws.Rows.Item(6).RowHeight = 105
$ws.Range("D6").Value = "USE mydatabase;`nCREATE TABLE mytable`n(`n name VARCHAR(100),`n    age INT`n);`n"
$ws.Range("D6").WrapText = $true

# --- Row 7: C7 "Drop " -> "To delete database"; D7 gets new "DROP DATABASE mydabase" ---
$ws.Range("C7").Value = "To delete database"
$ws.Range("D7").Value = "DROP DATABASE mydabase"

# --- Row 8: C8 and D8 get new content ---
$ws.Range("C8").Value = "To delete table"
$ws.Range("D8").Value = "DROP TABLE mytable"

# --- Sheet view: scroll so row 2 is at top, and select D4 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("D4").Select()
